$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before C ("level" -> D), creating space for the new "Name" column
$ws.Range("C1").EntireColumn.Insert()

# Insert a new column before the (now shifted) "KR phong" column (G), creating space for "krId"
$ws.Range("G1").EntireColumn.Insert()

# --- Header row ---
$ws.Range("B1").Value = "employeeId"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "level"
$ws.Range("E1").Value = "teamName"
$ws.Range("F1").Value = "Loại"
$ws.Range("G1").Value = "krId"
$ws.Range("H1").Value = "KR phòng"
$ws.Range("I1").Value = "KR team"
$ws.Range("J1").Value = "KR cá nhân"
$ws.Range("K1").Value = "Công thức tính"
$ws.Range("L1").Value = "Nguồn dữ liệu"
$ws.Range("M1").Value = "Định kỳ tính"
$ws.Range("N1").Value = "Đơn vị tính"
$ws.Range("O1").Value = "Điều kiện"
$ws.Range("P1").Value = "Norm"
$ws.Range("Q1").Value = "% Trọng số chỉ tiêu"
$ws.Range("R1").Value = "Kết quả"
$ws.Range("S1").Value = "Tỷ lệ"
$ws.Range("T1").Value = "Tổng thời gian dự kiến/ ước tính công việc (giờ)"
$ws.Range("U1").Value = "Tổng thời gian thực hiện công việc thực tế (giờ)"
$ws.Range("V1").Value = "Note"

# --- Row 2 ---
$ws.Range("A2").Value = 10
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = "pham"
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = "python"
$ws.Range("F2").Value = "KPI"
$ws.Range("G2").Value = 3
$ws.Range("H2").Value = "test kpi 4 1"
$ws.Range("I2").Value = "test kpi 4 1"
$ws.Range("J2").Value = "test kpi 4 1"
$ws.Range("K2").Value = "Báo cáo được CBQL confirm"
$ws.Range("L2").Value = "email"
$ws.Range("M2").Value = "Tháng"
$ws.Range("N2").Value = "%"
$ws.Range("O2").Value = "'="
$ws.Range("O2").Style = "Normal"
$ws.Range("P2").Value = 100
$ws.Range("Q2").Value = 80
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
# T2/U2 ("168"/"168") are left untouched, carried over as text from the original R2/S2 cells
# V2 (Note) is left blank/empty, carried over from the original empty T2 cell

# --- Row 3 ---
$ws.Range("A3").Value = 11
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = "pham"
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = "python"
$ws.Range("F3").Value = "KPI"
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = "test kpi 4 1"
$ws.Range("I3").Value = "test kpi 4 1"
$ws.Range("J3").Value = "test kpi 4 1"
$ws.Range("K3").Value = "Báo cáo được CBQL confirm"
$ws.Range("L3").Value = "email"
$ws.Range("M3").Value = "Tháng"
$ws.Range("N3").Value = "%"
$ws.Range("O3").Value = "'="
$ws.Range("O3").Style = "Normal"
$ws.Range("P3").Value = 100
$ws.Range("Q3").Value = 20
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0
# T3/U3 ("168"/"168") are left untouched, carried over as text from the original R3/S3 cells
# V3 (Note) is left blank/empty, carried over from the original empty T3 cell
